$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin -> Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "73.033.56"
$ws.Range("E2").Value = "  +2.17%  "

# Row 3: Ethereum -> Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "4.005.53"
$ws.Range("E3").Value = "  +0.72%  "

# Row 4: TetherUSD -> TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.18%  "

# Row 5: BNB -> BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "627.79"
$ws.Range("E5").Value = "  +16.67%  "

# Row 6: Solana -> Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "163.62"
$ws.Range("E6").Value = "  +8.47%  "

# Row 7: XRP -> XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.691"

# Row 8: USDC -> USDC
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.998"
$ws.Range("E8").Value = "  -0.09%  "

# Row 9: Cardano -> Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.765"
$ws.Range("E9").Value = "  +2.22%  "

# Row 10: Dogecoin -> Dogecoin
$ws.Range("E10").Value = "  +0.47%  "

# Row 11: Avalanche -> Avalanche
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.97"
$ws.Range("E11").Value = "  -0.89%  "

# Row 12: ShibaInu -> ShibaInu
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000320"
$ws.Range("E12").Value = "  -0.60%  "

# Row 13: Polkadot -> Polkadot
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.27"
$ws.Range("E13").Value = "  +5.05%  "

# Row 14: WrappedliquidstakedEther2.0 -> WrappedliquidstakedEther2.0
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.639.27"
$ws.Range("E14").Value = "  +0.66%  "

# Row 15: WrappedEther -> WrappedEther
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.003.07"
$ws.Range("E15").Value = "  +0.72%  "

# Row 16: Polygon -> Polygon
$ws.Range("E16").Value = "  +7.55%  "

# Row 17: Uniswap -> Uniswap
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.21"
$ws.Range("E17").Value = "  +1.00%  "

# Row 18: Chainlink -> Chainlink
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.76"
$ws.Range("E18").Value = "  +0.54%  "

# Row 19: TRON -> TRON
$ws.Range("E19").Value = "  +0.40%  "

# Row 20: WrappedBTC -> WrappedBTC
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.727.36"
$ws.Range("E20").Value = "  +1.95%  "

# Row 21: BitcoinCash -> BitcoinCash
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "442.43"
$ws.Range("E21").Value = "  +1.88%  "

# Row 22: PancakeSwap -> PancakeSwap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.15"
$ws.Range("E22").Value = "  +21.41%  "

# Row 23: Litecoin -> Litecoin
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "97.29"
$ws.Range("E23").Value = "  -0.21%  "

# Row 24: ImmutableX -> ImmutableX
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.47"
$ws.Range("E24").Value = "  -3.73%  "

# Row 25: InternetComputer(DFINITY) -> InternetComputer(DFINITY)
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.49"
$ws.Range("E25").Value = "  -1.24%  "

# Row 26: Toncoin -> Toncoin
$ws.Range("E26").Value = "  +5.37%  "

# Row 27: RenderToken -> RenderToken
$ws.Range("E27").Value = "  -0.14%  "

# Row 28: Filecoin -> Filecoin
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.64"
$ws.Range("E28").Value = "  -2.30%  "

# Row 29: LEO -> LEO
$ws.Range("E29").Value = "  +1.10%  "

# Row 30: EthereumClassic -> EthereumClassic
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.59"
$ws.Range("E30").Value = "  -0.71%  "

# Row 31: NEARProtocol -> NEARProtocol
$ws.Range("E31").Value = "  -2.87%  "

# Row 32: Cosmos -> Cosmos
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "14.09"
$ws.Range("E32").Value = "  +5.05%  "

# Row 33: Hedera -> Hedera
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.133"
$ws.Range("E33").Value = "  -0.35%  "

# Row 34: OKB -> OKB
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "72.25"
$ws.Range("E34").Value = "  +9.70%  "

# Row 35: InjectiveProtocol -> InjectiveProtocol
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "48.39"
$ws.Range("E35").Value = "  -6.73%  "

# Row 36: Bittensor -> Bittensor
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "639.94"
$ws.Range("E36").Value = "  -5.69%  "

# Row 37: PEPE -> PEPE
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0902"
$ws.Range("E37").Value = "  +9.05%  "

# Row 38: TheGraph -> TheGraph
$ws.Range("E38").Value = "  -1.30%  "

# Row 39: Kaspa -> Kaspa
$ws.Range("E39").Value = "  -0.93%  "

# Row 40: ThetaToken -> WEMIXToken
$ws.Range("B40").Value = "WEMIXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("E40").Value = "  +4.11%  "

# Row 41: WEMIXToken -> ThetaToken
$ws.Range("B41").Value = "ThetaToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.37"
$ws.Range("E41").Value = "  -1.38%  "

# Row 42: Dai -> Dai
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.998"
$ws.Range("E42").Value = "  -0.14%  "

# Row 43: FirstDigitalUSD -> FirstDigitalUSD
$ws.Range("E43").Value = "  +0.28%  "

# Row 44: VeChain -> VeChain
$ws.Range("E44").Value = "  +1.07%  "

# Row 45: THORChain -> THORChain
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.72"
$ws.Range("E45").Value = "  +2.47%  "

# Row 46: Stellar -> Stellar
$ws.Range("E46").Value = "  +1.06%  "

# Row 47: Fetch.AI -> Fetch.AI
$ws.Range("E47").Value = "  -0.68%  "

# Row 48: ApeXProtocol -> ApeXProtocol
$ws.Range("E48").Value = "  +2.56%  "

# Row 49: Maker -> Stacks
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.09"
$ws.Range("E49").Value = "  +2.00%  "

# Row 50: Stacks -> Maker
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.903.35"
$ws.Range("E50").Value = "  +10.28%  "

# Row 51: LidoDAOToken -> LidoDAOToken
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.44"
$ws.Range("E51").Value = "  +4.08%  "
